$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J: DATE_TYPE_CODE changes from "001" to "004" (keep as text, like the original code)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "004"

# REPORT_DATE changes from 2019-12-31 to 2020-09-30
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Numeric values update to the new reporting period figures
$ws.Range("O2").Value = 22981720.11
$ws.Range("P2").Value = 11.2710424161
$ws.Range("Q2").Value = 317072528.57
$ws.Range("R2").Value = 155.503500233
$ws.Range("S2").Value = 69435306.92
$ws.Range("T2").Value = 34.0535123447
$ws.Range("U2").Value = -190846024.58
$ws.Range("V2").Value = -93.59759094109999
$ws.Range("W2").Value = 2254928.47
$ws.Range("X2").Value = 1.1058960908
$ws.Range("Y2").Value = 160608001.47
$ws.Range("Z2").Value = 78.76780276949999
$ws.Range("AA2").Value = -36036268.89
$ws.Range("AB2").Value = -17.673451475
$ws.Range("AC2").Value = -203900573.36
$ws.Range("AD2").Value = -159.3269056566
